$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: CheckIn
$ws.Cells.Item(2,1).Value2 = "CheckIn"
$ws.Cells.Item(2,3).Value2 = 2
$ws.Cells.Item(2,4).Value2 = '["Business", "Speakers", "Logistics"]'

# Row 3: Auditorio
$ws.Cells.Item(3,1).Value2 = "Auditorio"
$ws.Cells.Item(3,3).Value2 = 6
$ws.Cells.Item(3,4).Value2 = '["Logistics"]'

# Row 4: Almocos
$ws.Cells.Item(4,1).Value2 = "Almocos"
$ws.Cells.Item(4,3).Value2 = 6
$ws.Cells.Item(4,4).Value2 = '["Business", "Logistics"]'

# Row 5: Divulgacao
$ws.Cells.Item(5,1).Value2 = "Divulgacao"
$ws.Cells.Item(5,3).Value2 = 3
$ws.Cells.Item(5,4).Value2 = '["Marketing", "Volunteer"]'

# Row 6: Workshops
$ws.Cells.Item(6,1).Value2 = "Workshops"
$ws.Cells.Item(6,3).Value2 = 2
$ws.Cells.Item(6,4).Value2 = '[]'

# Widen column B (manual resize by author; no longer "best fit")
$ws.Columns.Item(2).ColumnWidth = 37
